$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 9 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(9, 1).Value = "Sollicitatie salesfunctie"
$logs.Cells.Item(9, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(9, 3).Value = "Hierbij mijn sollicitatie voor de salesfunctie. CV in bijlage."
$logs.Cells.Item(9, 4).Value = "Sollicitatie / Vacature"
$logs.Cells.Item(9, 5).Value = "Beste,`nBedankt voor je interesse in onze openstaande salesfunctie. Ik heb je sollicitatie ontvangen en zal deze zo spoedig mogelijk bekijken. Mocht ik nog aanvullende informatie nodig hebben, dan neem ik contact met je op.`nMet vriendelijke groet,`n[Naam Bedrijf]"
$logs.Cells.Item(9, 6).Value = "2025-06-24 19:51:02"
$logs.Cells.Item(9, 7).Value = "Ja"

# Setting multi-line text auto-expands the row height; restore the
# default (matches the other data rows, none of which have a custom height).
$logs.Rows.Item(9).RowHeight = 15

# --- Extend conditional formatting ranges to cover the new row ---
$dFc = $logs.Range("D2:D8").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D9"))
}

$gFc = $logs.Range("G2:G8").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G9"))
}

# --- Sheet "Dashboard": Sollicitatie / Vacature count goes to 2, swaps order with Factuur / Administratie ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(3, 2).Value = 2

$dash.Cells.Item(4, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(4, 2).Value = 2
